$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 412, which shifts rows 412:517 down to 413:518
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the new record's values.
# Columns A,B,C,E,F,G,H,I,R are identical to the rest of this data block.
$ws.Cells.Item(412, 1).Value = 10
$ws.Cells.Item(412, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(412, 3).Value = "La Araucanía"
$ws.Cells.Item(412, 4).Value = 45204
$ws.Cells.Item(412, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(412, 5).Value = 9
$ws.Cells.Item(412, 6).Value = 100112001
$ws.Cells.Item(412, 7).Value = "Berenjena"
$ws.Cells.Item(412, 8).Value = "Sin especificar"
$ws.Cells.Item(412, 9).Value = "Primera"
$ws.Cells.Item(412, 10).Value = 35
$ws.Cells.Item(412, 11).Value = 12000
$ws.Cells.Item(412, 12).Value = 12000
$ws.Cells.Item(412, 13).Value = 12000
$ws.Cells.Item(412, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(412, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(412, 16).Value = 300
$ws.Cells.Item(412, 17).Value = 40
$ws.Cells.Item(412, 18).Value = "Hortaliza"
